$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 0.1 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text (shared across all test-case blocks: B8, B15, B22, B30, B37, B45, B53)
$precondition = "O usuário devidamente autenticado e na tela de listagem de empenhos."
$ws.Range("B8").Value = $precondition
$ws.Range("B15").Value = $precondition
$ws.Range("B22").Value = $precondition
$ws.Range("B30").Value = $precondition
$ws.Range("B37").Value = $precondition
$ws.Range("B45").Value = $precondition
$ws.Range("B53").Value = $precondition

# TC1 - expected result (D10): add trailing period
$ws.Range("D10").Value = "SYSTEM Recupera e exibe todos os detalhes (dados) da solicitação para o usuário; e Apresenta a tela de Detalhar Diárias."

# TC3 - expected result (D24): wording + accent fixes
$ws.Range("D24").Value = "SYSTEM Exibe a lista de solicitações aguardando serem empenhadas, de todos os servidores, ordenado pelo número da diária em ordem crescente."

# TC4 - expected result (D32): add trailing period
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Registrar Empenho."

# TC5 - action step text (shared across B39, B47, B55): capitalize "Filtra" + trailing period
$tc5Action = "Chefe/Beneficiário Filtra a listagem por registros cujos beneficiários não possuem número do credor."
$ws.Range("B39").Value = $tc5Action
$ws.Range("B47").Value = $tc5Action
$ws.Range("B55").Value = $tc5Action
